$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EXPD")

# Insert a new column before D (shifts old D:K -> E:L)
$ws.Range("D:D").Insert()

# Copy number formatting/style from new column E (the old column D) into the
# newly inserted column D so each row gets the right style (date vs number)
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest fiscal-year figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 8138400
$ws.Range("D9").Value = 5518000
$ws.Range("D10").Value = 2620400
$ws.Range("D15").Value = 54000
$ws.Range("D17").Value = 7341800
$ws.Range("D18").Value = 796600
$ws.Range("D20").Value = 21800
$ws.Range("D21").Value = 872300
$ws.Range("D22").Value = "NA"
$ws.Range("D23").Value = 818300
$ws.Range("D24").Value = 179500
$ws.Range("D26").Value = 638800
$ws.Range("D27").Value = 637200
$ws.Range("D29").Value = -19000
$ws.Range("D32").Value = -21800
$ws.Range("D33").Value = 618200
$ws.Range("D35").Value = 618200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 923700
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 1581500
$ws.Range("D45").Value = 229600
$ws.Range("D46").Value = 2734800
$ws.Range("D48").Value = 504100
$ws.Range("D49").Value = 7900
$ws.Range("D52").Value = 67700
$ws.Range("D54").Value = 3314600
$ws.Range("D57").Value = 902300
$ws.Range("D59").Value = 424600
$ws.Range("D60").Value = 1326800
$ws.Range("D62").Value = 0
$ws.Range("D66").Value = 1327700
$ws.Range("D72").Value = 2088700
$ws.Range("D76").Value = 1986800
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 618200
$ws.Range("D83").Value = 54000
$ws.Range("D89").Value = 572800
$ws.Range("D91").Value = -47500
$ws.Range("D94").Value = -48400
$ws.Range("D96").Value = -156800
$ws.Range("D100").Value = -627700
$ws.Range("D101").Value = -24000
$ws.Range("D102").Value = -127400

# A few small historical corrections that came along with this edit
$ws.Range("F89").Value = 529500
$ws.Range("F100").Value = -298800
$ws.Range("K83").Value = 37900
